$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.542.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.844.55'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -1.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.10'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4640'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3847'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.05'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07899'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9972'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.50'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.849.72'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.966'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.124'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06678'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001035'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.14'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.529.61'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.390'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.87'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.306'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.11'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.057.56'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.51'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.116'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.412'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.90'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9753'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.00%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.590'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.301'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.334'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06039'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02228'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.283'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.181'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5895'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.36'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.14%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5589'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.15'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.912'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06695'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.47'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.050'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.42%  '
